$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.34533969975764
$ws.Range("C2").Value = 8.16920413881499
$ws.Range("E2").Value = 19.85142983026525
$ws.Range("F2").Value = 39.95938947265507
$ws.Range("G2").Value = 32.14806665429052
$ws.Range("H2").Value = 14.67659151876387
$ws.Range("J2").Value = 7.924973918464391
$ws.Range("M2").Value = 18.99327997478041
$ws.Range("N2").Value = 17.53516521928421
$ws.Range("B3").Value = 12.80444162561037
$ws.Range("C3").Value = 7.687040072786671
$ws.Range("E3").Value = 19.84197809585668
$ws.Range("F3").Value = 39.83086723860232
$ws.Range("G3").Value = 31.87143442498389
$ws.Range("H3").Value = 14.70279423431628
$ws.Range("J3").Value = 7.949473396150185
$ws.Range("M3").Value = 18.83316733293717
$ws.Range("N3").Value = 17.60699527740569
$ws.Range("B4").Value = 12.46334788739799
$ws.Range("C4").Value = 7.373721007158663
$ws.Range("E4").Value = 19.83897594606293
$ws.Range("F4").Value = 39.76312424993311
$ws.Range("G4").Value = 31.71467284520067
$ws.Range("H4").Value = 14.72260957639005
$ws.Range("J4").Value = 7.965272384687556
$ws.Range("M4").Value = 18.73796815373505
$ws.Range("N4").Value = 17.65305572237402
$ws.Range("B5").Value = 12.32231687823949
$ws.Range("C5").Value = 7.241726594456999
$ws.Range("E5").Value = 19.83845915653813
$ws.Range("F5").Value = 39.73834382026966
$ws.Range("G5").Value = 31.65414842924466
$ws.Range("H5").Value = 14.7316181621001
$ws.Range("J5").Value = 7.971901304004829
$ws.Range("M5").Value = 18.69999187506637
$ws.Range("N5").Value = 17.67231945874966
$ws.Range("B6").Value = 12.29878312040093
$ws.Range("C6").Value = 7.219549433352671
$ws.Range("E6").Value = 19.83841608072136
$ws.Range("F6").Value = 39.73440010777152
$ws.Range("G6").Value = 31.64430289972205
$ws.Range("H6").Value = 14.73317031562486
$ws.Range("J6").Value = 7.973013565144783
$ws.Range("M6").Value = 18.69373637926844
$ws.Range("N6").Value = 17.67554805915985
$ws.Range("B7").Value = 12.46145380673076
$ws.Range("C7").Value = 7.371958309834289
$ws.Range("E7").Value = 19.83896611270746
$ws.Range("F7").Value = 39.76277859344034
$ws.Range("G7").Value = 31.71384291891957
$ws.Range("H7").Value = 14.72272729376443
$ws.Range("J7").Value = 7.965361011742106
$ws.Range("M7").Value = 18.73745263328807
$ws.Range("N7").Value = 17.65331351846893
$ws.Range("B8").Value = 13.16082938964994
$ws.Range("C8").Value = 8.006552226527946
$ws.Range("E8").Value = 19.84759087761158
$ws.Range("F8").Value = 39.91276765339988
$ws.Range("G8").Value = 32.05000836497974
$ws.Range("H8").Value = 14.68485073475112
$ws.Range("J8").Value = 7.933264717740326
$ws.Range("M8").Value = 18.93744954368987
$ws.Range("N8").Value = 17.55952721239984
$ws.Range("B9").Value = 14.45240050068983
$ws.Range("C9").Value = 9.113139198998567
$ws.Range("E9").Value = 19.88662427310974
$ws.Range("F9").Value = 40.29463406277179
$ws.Range("G9").Value = 32.80983418519806
$ws.Range("H9").Value = 14.64028734996751
$ws.Range("J9").Value = 7.876298952589872
$ws.Range("M9").Value = 19.35267883329326
$ws.Range("N9").Value = 17.391055200576
$ws.Range("B10").Value = 15.34267257545531
$ws.Range("C10").Value = 9.841435358109653
$ws.Range("E10").Value = 19.92863092077555
$ws.Range("F10").Value = 40.62725127844843
$ws.Range("G10").Value = 33.42450991071703
$ws.Range("H10").Value = 14.62583613402054
$ws.Range("J10").Value = 7.838053109635728
$ws.Range("M10").Value = 19.66950988145984
$ws.Range("N10").Value = 17.27657744362229
$ws.Range("B11").Value = 15.73327515506662
$ws.Range("C11").Value = 10.1543127308435
$ws.Range("E11").Value = 19.95059591049917
$ws.Range("F11").Value = 40.78949641777596
$ws.Range("G11").Value = 33.71522459070309
$ws.Range("H11").Value = 14.62326391841995
$ws.Range("J11").Value = 7.821429723947018
$ws.Range("M11").Value = 19.81569643860344
$ws.Range("N11").Value = 17.22649290389953
$ws.Range("B12").Value = 15.87900324739635
$ws.Range("C12").Value = 10.27014088301223
$ws.Range("E12").Value = 19.95932029072553
$ws.Range("F12").Value = 40.85247031873255
$ws.Range("G12").Value = 33.82680137318442
$ws.Range("H12").Value = 14.62286719665859
$ws.Range("J12").Value = 7.815245739299795
$ws.Range("M12").Value = 19.87130663993563
$ws.Range("N12").Value = 17.20781178495174
$ws.Range("B13").Value = 15.84771692688334
$ws.Range("C13").Value = 10.24531315432657
$ws.Range("E13").Value = 19.95742331689732
$ws.Range("F13").Value = 40.83884011123383
$ws.Range("G13").Value = 33.8027067974899
$ws.Range("H13").Value = 14.62292693554154
$ws.Range("J13").Value = 7.816572645113015
$ws.Range("M13").Value = 19.85931945051038
$ws.Range("N13").Value = 17.21182245708039
$ws.Range("B14").Value = 15.7453086809881
$ws.Range("C14").Value = 10.1638951278271
$ws.Range("E14").Value = 19.95130554310884
$ws.Range("F14").Value = 40.79464677573812
$ws.Range("G14").Value = 33.72437484858343
$ws.Range("H14").Value = 14.62321970015219
$ws.Range("J14").Value = 7.820918743381107
$ws.Range("M14").Value = 19.82026671185334
$ws.Range("N14").Value = 17.22495029831763
$ws.Range("B15").Value = 15.68229301526738
$ws.Range("C15").Value = 10.11367889631882
$ws.Range("E15").Value = 19.94761107354592
$ws.Range("F15").Value = 40.76777584815186
$ws.Range("G15").Value = 33.6765850314051
$ws.Range("H15").Value = 14.62347426017774
$ws.Range("J15").Value = 7.823595286056728
$ws.Range("M15").Value = 19.7963773752416
$ws.Range("N15").Value = 17.23302852241485
$ws.Range("B16").Value = 15.31684599573169
$ws.Range("C16").Value = 9.820617050726248
$ws.Range("E16").Value = 19.92725259379165
$ws.Range("F16").Value = 40.61686510862329
$ws.Range("G16").Value = 33.40572530561256
$ws.Range("H16").Value = 14.62608492857957
$ws.Range("J16").Value = 7.839155036710266
$ws.Range("M16").Value = 19.65999395493194
$ws.Range("N16").Value = 17.27989052969317
$ws.Range("B17").Value = 15.08888350078759
$ws.Range("C17").Value = 9.636111062518605
$ws.Range("E17").Value = 19.91549207611666
$ws.Range("F17").Value = 40.52706157159137
$ws.Range("G17").Value = 33.24232915355119
$ws.Range("H17").Value = 14.62871288413101
$ws.Range("J17").Value = 7.848898544281227
$ws.Range("M17").Value = 19.57682349719861
$ws.Range("N17").Value = 17.30914787826077
$ws.Range("B18").Value = 14.95641715386186
$ws.Range("C18").Value = 9.528253460581288
$ws.Range("E18").Value = 19.90899680985514
$ws.Range("F18").Value = 40.47644072690063
$ws.Range("G18").Value = 33.14939898862414
$ws.Range("H18").Value = 14.63060098714542
$ws.Range("J18").Value = 7.854575712743792
$ws.Range("M18").Value = 19.52918242790133
$ws.Range("N18").Value = 17.32616352162635
$ws.Range("B19").Value = 14.91133855119821
$ws.Range("C19").Value = 9.491436902818347
$ws.Range("E19").Value = 19.90684395374671
$ws.Range("F19").Value = 40.45947967625011
$ws.Range("G19").Value = 33.11811811241991
$ws.Range("H19").Value = 14.63130487941058
$ws.Range("J19").Value = 7.856510450511835
$ws.Range("M19").Value = 19.51308703085084
$ws.Range("N19").Value = 17.33195699481353
$ws.Range("B20").Value = 15.11329093877626
$ws.Range("C20").Value = 9.655931589450569
$ws.Range("E20").Value = 19.91671618315118
$ws.Range("F20").Value = 40.53651478815596
$ws.Range("G20").Value = 33.25961493822236
$ws.Range("H20").Value = 14.62839414523679
$ws.Range("J20").Value = 7.847853783715113
$ws.Range("M20").Value = 19.58565712463114
$ws.Range("N20").Value = 17.30601398062343
$ws.Range("B21").Value = 15.77544859676328
$ws.Range("C21").Value = 10.18788154377016
$ws.Range("E21").Value = 19.95309147651916
$ws.Range("F21").Value = 40.80758607771278
$ws.Range("G21").Value = 33.74734327967366
$ws.Range("H21").Value = 14.62311802674083
$ws.Range("J21").Value = 7.819639182220628
$ws.Range("M21").Value = 19.83173094152654
$ws.Range("N21").Value = 17.22108661644164
$ws.Range("B22").Value = 16.1954307706053
$ws.Range("C22").Value = 10.52008627295708
$ws.Range("E22").Value = 19.97923343313589
$ws.Range("F22").Value = 40.99367569529672
$ws.Range("G22").Value = 34.07473221929477
$ws.Range("H22").Value = 14.62303543732668
$ws.Range("J22").Value = 7.801845674665553
$ws.Range("M22").Value = 19.99400611172129
$ws.Range("N22").Value = 17.16724101042087
$ws.Range("B23").Value = 15.97248148284082
$ws.Range("C23").Value = 10.34419644753171
$ws.Range("E23").Value = 19.96506563891949
$ws.Range("F23").Value = 40.89355227971416
$ws.Range("G23").Value = 33.89924476883189
$ws.Range("H23").Value = 14.62277103971425
$ws.Range("J23").Value = 7.811283422646127
$ws.Range("M23").Value = 19.90727835308645
$ws.Range("N23").Value = 17.19582813129635
$ws.Range("B24").Value = 15.10226071167477
$ws.Range("C24").Value = 9.646976281689405
$ws.Range("E24").Value = 19.91616193548138
$ws.Range("F24").Value = 40.53223784506622
$ws.Range("G24").Value = 33.25179688570909
$ws.Range("H24").Value = 14.628537072104
$ws.Range("J24").Value = 7.848325884763612
$ws.Range("M24").Value = 19.58166289408301
$ws.Range("N24").Value = 17.30743020756254
$ws.Range("B25").Value = 14.11266205701796
$ws.Range("C25").Value = 8.828625168294273
$ws.Range("E25").Value = 19.87371153108097
$ws.Range("F25").Value = 40.18206968291963
$ws.Range("G25").Value = 32.59399072758377
$ws.Range("H25").Value = 14.64914247539237
$ws.Range("J25").Value = 7.891073738830993
$ws.Range("M25").Value = 19.23811402195234
$ws.Range("N25").Value = 17.43499009331327
